# Applies the diff:
#  - moves the "_GoBack" bookmark from right after the second
#    "Description/Overview:" run to span from the start of the second
#    "Iteration: l1" paragraph through the end of the document's last
#    paragraph ("*006 - Purchase Order Approval Protocol")
#  - changes the second "Iteration: l1" paragraph's text to
#    "Iteration: E1" (split as "Iteration: " / "E" / "1" in the source,
#    kept bold throughout)

$d = $word.ActiveDocument

# Locate the target paragraphs *before* editing anything, since paragraph
# indices stay stable across the edits we are about to make (no paragraphs
# are added or removed).
$paraCount = $d.Paragraphs.Count
$iterationParaIndex = 0
$matchCount = 0
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    $candidateText = $candidate.Range.Text.TrimEnd([char]13)
    if ($candidateText -eq "Iteration: l1") {
        $matchCount = $matchCount + 1
        if ($matchCount -eq 2) {
            $iterationParaIndex = $i
        }
    }
}

# Remove the existing "_GoBack" bookmark (currently sitting right after the
# second "Description/Overview:" run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Split "Iteration: l1" into "Iteration: " + "E" + "1", all bold.
$iterationPara = $d.Paragraphs.Item($iterationParaIndex)
$fullParaRange = $iterationPara.Range
$textRange = $d.Range($fullParaRange.Start, $fullParaRange.End - 1)
$iterationParaStart = $textRange.Start

$textRange.Text = "Iteration: "

$insertE = $d.Range($textRange.End, $textRange.End)
$insertE.InsertAfter("E")

$insert1 = $d.Range($textRange.End + 1, $textRange.End + 1)
$insert1.InsertAfter("1")

# Re-add the "_GoBack" bookmark spanning from the start of the
# "Iteration: ..." paragraph through the end of the last paragraph's text
# (i.e. right after "*006 - Purchase Order Approval Protocol", before its
# trailing paragraph mark).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkEndPos = $lastPara.Range.End - 1
$bookmarkRange = $d.Range($iterationParaStart, $bookmarkEndPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
